$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric results (rows 3-8) ---
$ws.Range("A3").Value = 15.333453060715755
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 15.333453060715755
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 15.524891059342343
$ws.Range("F3").Value = 0.10793285726290715
$ws.Range("G3").Value = 15.333453060716749
$ws.Range("H3").Value = 0.0000000000022267341529379978

$ws.Range("A4").Value = 199.6030061177791
$ws.Range("B4").Value = 0.00086100158977057181
$ws.Range("C4").Value = 199.60420634922809
$ws.Range("D4").Value = 0.00088732851388176597
$ws.Range("E4").Value = 199.6046031746032
$ws.Range("F4").Value = 0.000000000000028421709430404007
$ws.Range("G4").Value = 199.60301587272761
$ws.Range("H4").Value = 0.00088732790601123363

$ws.Range("A5").Value = 127291.00401606425
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 127291.00401606425
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 127291.0040160887
$ws.Range("F5").Value = 0.000000051357248764017916
$ws.Range("G5").Value = 127291.00401606425
$ws.Range("H5").Value = 0

$ws.Range("A6").Value = 4535.8819444444443
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 4535.8819444444443
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 4535.8838383252532
$ws.Range("F6").Value = 0.000000070180116077090753
$ws.Range("G6").Value = 4535.8819444444443
$ws.Range("H6").Value = 0

$ws.Range("A7").Value = 502053.26589279092
$ws.Range("B7").Value = 0.050857770806125913
$ws.Range("C7").Value = 502056.5035258549
$ws.Range("D7").Value = 1.8213603576168393
$ws.Range("E7").Value = 502058.49627535336
$ws.Range("F7").Value = 0.97348043596611722
$ws.Range("G7").Value = 502053.24884617608
$ws.Range("H7").Value = 0.015735821426437174

$ws.Range("A8").Value = 274.04065217391292
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 274.04065217391292
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 274.04065217391292
$ws.Range("F8").Value = 0.000000000000028421709430404007
$ws.Range("G8").Value = 274.04065217391292
$ws.Range("H8").Value = 0

# B/D/F/H (STD columns) now share the same (uncoloured) look as A/C/E/G
for ($r = 3; $r -le 8; $r++) {
    $ws.Range("A$r").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)
    $ws.Range("C$r").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)
    $ws.Range("E$r").Copy()
    $ws.Range("F$r").PasteSpecial(-4122)
    $ws.Range("G$r").Copy()
    $ws.Range("H$r").PasteSpecial(-4122)
}

# --- Rows 9-10 no longer hold a results set: clear the values, keep formatting ---
$ws.Range("A9:H10").ClearContents()

# --- Selection moves from C6 to C9 ---
$ws.Range("C9").Select()
